# Updated the Spread Sheet
# Adds a new "Assignment_7" row (row 8) to the Assignments tracker, mirroring
# the formatting/layout of the preceding row (Assignment_6, row 7), including
# its hyperlink in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the look (borders, fonts, alignment, etc.) of the last existing data
# row (row 7) down onto the new row 8 before we populate it.
$ws.Range("A7:C7").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row's data.
$ws.Range("A8").Value = "Assignment_7"
$ws.Range("B8").Value = "https://github.com/Vasanth30e/Assignments_Phase2/tree/master/Assignment_7"
$ws.Hyperlinks.Add($ws.Range("B8"), "https://github.com/Vasanth30e/Assignments_Phase2/tree/master/Assignment_7")
$ws.Range("C8").Value = 45159

# Match the taller row height used for the new entry.
$ws.Rows("8").RowHeight = 33.75

# Reflect the author's final cursor position in the saved view.
$ws.Range("E9").Select()
